$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '29.462.85'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +0.61%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.877.97'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +1.00%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.7134'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +1.34%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '242.16'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +1.57%  '
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.07847'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -1.80%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3115'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +2.86%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '25.17'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +7.00%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.08265'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +0.83%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.7310'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +3.26%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.275'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +1.38%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '1.854.07'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -1.10%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '91.20'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +1.62%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '29.454.79'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +0.19%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '5.923'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +1.38%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '247.56'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +3.90%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000007879'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -0.44%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.29'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -0.17%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.9996'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -0.17%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '7.955'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +6.28%  '
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.1589'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +10.32%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '163.75'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +0.45%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.016'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +1.33%  '
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +0.95%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.364'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -4.04%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.492'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +0.96%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.373'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -0.16%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.132'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +2.52%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.05316'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +2.30%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.935'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +0.49%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.201'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +3.24%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.7238'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +0.89%  '
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -0.13%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.01869'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +0.82%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.263.33'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +9.29%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.738'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +0.34%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.9109'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -3.64%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '74.06'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +4.63%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '6.143'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +2.20%  '
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '103.43'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +0.40%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.5330'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +0.53%  '
$ws.Range('B46').NumberFormat = "@"
$ws.Range('B46').Value = 'RocketPoolETH'
$ws.Range('C46').NumberFormat = "@"
$ws.Range('C46').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.011.91'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -1.27%  '
$ws.Range('B47').NumberFormat = "@"
$ws.Range('B47').Value = 'SynthetixNetwork'
$ws.Range('C47').NumberFormat = "@"
$ws.Range('C47').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.961'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +14.12%  '
$ws.Range('B48').NumberFormat = "@"
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').NumberFormat = "@"
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.773'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +0.41%  '
$ws.Range('B49').NumberFormat = "@"
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').NumberFormat = "@"
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.00000000120'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -0.32%  '
$ws.Range('B50').NumberFormat = "@"
$ws.Range('B50').Value = 'TheSandbox'
$ws.Range('C50').NumberFormat = "@"
$ws.Range('C50').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.4326'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +1.26%  '
$ws.Range('B51').NumberFormat = "@"
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').NumberFormat = "@"
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '9.250'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +0.84%  '
